$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.248.06"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "1.645.22"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.50"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.09"
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "1.875.44"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "1.644.71"
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.546"
$ws.Range("E15").Value = "  +3.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.45"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").Value = "27.226.50"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.96"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.92"
$ws.Range("E21").Value = "  +4.90%  "
$ws.Range("E22").Value = "  +7.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.41"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.95"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.54"
$ws.Range("E26").Value = "  +1.84%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.04"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").Value = "1.276.78"
$ws.Range("E35").Value = "  +2.50%  "
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.863"
$ws.Range("E38").Value = "  +4.08%  "
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  +6.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.32"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").Value = "1.785.33"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.96"
$ws.Range("E45").Value = "  +1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.03"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0516"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.70"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("E51").Value = "  +0.67%  "
